$wb = $excel.ActiveWorkbook

# --- "Toggles input del" sheet (sheet4): drop the 4th observation (row 5),
#     and change the remaining three toggle counts (experiment now uses only
#     one delayed input) -------------------------------------------------
$wsToggles = $wb.Worksheets.Item("Toggles input del")

# New toggle counts for the remaining 3 rows.
$wsToggles.Range("B2").Value = 96
$wsToggles.Range("B3").Value = 128
$wsToggles.Range("B4").Value = 32

# Remove the now-unused 4th observation (row 5) entirely so the sheet's
# used range / dimension shrinks from A1:B5 to A1:B4.
$wsToggles.Range("A5:B5").Delete()

# --- Chart embedded on that sheet references the toggle column; point it
#     at the new, smaller range -----------------------------------------
$chart = $wsToggles.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,,'Toggles input del'!`$B`$2:`$B`$4,1)"

# --- "Correlation matrix" sheet (sheet1): update the two correlations that
#     depend on the "gate+inputs delay" / toggle data that just changed ---
$wsCorr = $wb.Worksheets.Item("Correlation matrix")
$wsCorr.Range("E4").Value = 0.3779644730092272
$wsCorr.Range("G4").Value = 0.29277002188456
